$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# Row 6 (Item ID 4564)
$ws_ALC.Range("H6").Value = 22.666666
$ws_ALC.Range("I6").Value = 22.666666
$ws_ALC.Range("K6").Value = 67.99999800000001
$ws_ALC.Range("M6").Value = 44.00000199999999

# Row 15 (Item ID 44146)
$ws_ALC.Range("H15").Value = 501
$ws_ALC.Range("I15").Value = 501
$ws_ALC.Range("K15").Value = 1503
$ws_ALC.Range("M15").Value = -1334

# Row 33 (Item ID 5512)
$ws_ALC.Range("H33").Value = 139.07692
$ws_ALC.Range("I33").Value = 145.9
$ws_ALC.Range("K33").Value = 145.9
$ws_ALC.Range("M33").Value = 83.09999999999999

# Row 118 (Item ID 27958)
$ws_ALC.Range("H118").Value = 626.3333
$ws_ALC.Range("I118").Value = 626.3333
$ws_ALC.Range("K118").Value = 1878.9999
$ws_ALC.Range("M118").Value = -221.9999

# Row 125 (Item ID 36228)
$ws_ALC.Range("H125").Value = 1647.5
$ws_ALC.Range("I125").Value = 1125
$ws_ALC.Range("J125").Value = 1908.75
$ws_ALC.Range("K125").Value = 10125
$ws_ALC.Range("L125").Value = 17178.75
$ws_ALC.Range("M125").Value = -7665
$ws_ALC.Range("N125").Value = -22098.75

# --- ARM ---
# Row 2 (Item ID 27713)
$ws_ARM.Range("H2").Value = 730.5
$ws_ARM.Range("I2").Value = 730.5
$ws_ARM.Range("K2").Value = 730.5
$ws_ARM.Range("M2").Value = -617.5

# Row 4 (Item ID 5071)
$ws_ARM.Range("H4").Value = 439.8
$ws_ARM.Range("I4").Value = 300
$ws_ARM.Range("J4").Value = 999
$ws_ARM.Range("K4").Value = 300
$ws_ARM.Range("L4").Value = 999
$ws_ARM.Range("M4").Value = -184
$ws_ARM.Range("N4").Value = -1231

# Row 32 (Item ID 44147)
$ws_ARM.Range("H32").Value = 4251.8
$ws_ARM.Range("I32").Value = 4251.8
$ws_ARM.Range("K32").Value = 4251.8
$ws_ARM.Range("M32").Value = -3964.8

# Row 97 (Item ID 19941)
$ws_ARM.Range("H97").Value = 700.875
$ws_ARM.Range("I97").Value = 700.875
$ws_ARM.Range("K97").Value = 700.875
$ws_ARM.Range("M97").Value = -204.875

# Row 102 (Item ID 19945)
$ws_ARM.Range("H102").Value = 1166.6666
$ws_ARM.Range("I102").Value = 1250
$ws_ARM.Range("J102").Value = 1000
$ws_ARM.Range("K102").Value = 1250
$ws_ARM.Range("L102").Value = 1000
$ws_ARM.Range("M102").Value = 372
$ws_ARM.Range("N102").Value = -4244

# Row 110 (Item ID 27708)
$ws_ARM.Range("H110").Value = 999.6667
$ws_ARM.Range("I110").Value = 999.6667
$ws_ARM.Range("K110").Value = 999.6667
$ws_ARM.Range("M110").Value = 1045.3333

# Row 113 (Item ID 26002)
$ws_ARM.Range("H113").Value = 35000
$ws_ARM.Range("J113").Value = 35000
$ws_ARM.Range("L113").Value = 35000
$ws_ARM.Range("N113").Value = -43678

# Row 116 (Item ID 27713)
$ws_ARM.Range("H116").Value = 730.5
$ws_ARM.Range("I116").Value = 730.5
$ws_ARM.Range("K116").Value = 730.5
$ws_ARM.Range("M116").Value = 1563.5

# --- BSM ---
# Row 3 (Item ID 27713)
$ws_BSM.Range("H3").Value = 730.5
$ws_BSM.Range("I3").Value = 730.5
$ws_BSM.Range("K3").Value = 730.5
$ws_BSM.Range("M3").Value = -616.5

# Row 42 (Item ID 22903)
$ws_BSM.Range("H42").Value = 299999
$ws_BSM.Range("J42").Value = 299999
$ws_BSM.Range("L42").Value = 299999
$ws_BSM.Range("N42").Value = -300655

# Row 54 (Item ID 2376)
$ws_BSM.Range("H54").Value = 5148.1665
$ws_BSM.Range("I54").Value = 5148.1665
$ws_BSM.Range("J54").Value = 0
$ws_BSM.Range("K54").Value = 5148.1665
$ws_BSM.Range("L54").Value = 0
$ws_BSM.Range("M54").Value = -4664.1665
$ws_BSM.Range("N54").ClearContents()

# Row 82 (Item ID 11877)
$ws_BSM.Range("H82").Value = 5831
$ws_BSM.Range("I82").Value = 5831
$ws_BSM.Range("K82").Value = 5831
$ws_BSM.Range("M82").Value = -5448

# Row 85 (Item ID 11877)
$ws_BSM.Range("H85").Value = 5831
$ws_BSM.Range("I85").Value = 5831
$ws_BSM.Range("K85").Value = 5831
$ws_BSM.Range("M85").Value = -4505

# Row 86 (Item ID 12526)
$ws_BSM.Range("H86").Value = 2000
$ws_BSM.Range("I86").Value = 2000
$ws_BSM.Range("K86").Value = 2000
$ws_BSM.Range("M86").Value = -877

# Row 89 (Item ID 12526)
$ws_BSM.Range("H89").Value = 2000
$ws_BSM.Range("I89").Value = 2000
$ws_BSM.Range("K89").Value = 10000
$ws_BSM.Range("M89").Value = -4384

# Row 94 (Item ID 19939)
$ws_BSM.Range("H94").Value = 1323.25
$ws_BSM.Range("I94").Value = 1323.25
$ws_BSM.Range("K94").Value = 1323.25
$ws_BSM.Range("M94").Value = -872.25

# Row 105 (Item ID 19947)
$ws_BSM.Range("H105").Value = 2006.9
$ws_BSM.Range("I105").Value = 1952.2222
$ws_BSM.Range("K105").Value = 1952.2222
$ws_BSM.Range("M105").Value = -205.2221999999999

# Row 134 (Item ID 43998)
$ws_BSM.Range("H134").Value = 3420.4285
$ws_BSM.Range("I134").Value = 1388.7
$ws_BSM.Range("K134").Value = 4166.1
$ws_BSM.Range("M134").Value = -1631.1

# --- CRP ---
# Row 16 (Item ID 27691)
$ws_CRP.Range("H16").Value = 2556
$ws_CRP.Range("I16").Value = 491
$ws_CRP.Range("J16").Value = 3932.6667
$ws_CRP.Range("K16").Value = 491
$ws_CRP.Range("L16").Value = 3932.6667
$ws_CRP.Range("M16").Value = -204
$ws_CRP.Range("N16").Value = -4506.6667

# Row 39 (Item ID 1915)
$ws_CRP.Range("H39").Value = 5000
$ws_CRP.Range("I39").Value = 0
$ws_CRP.Range("J39").Value = 5000
$ws_CRP.Range("K39").Value = 0
$ws_CRP.Range("L39").Value = 5000
$ws_CRP.Range("M39").ClearContents()
$ws_CRP.Range("N39").Value = -5782

# Row 49 (Item ID 1915)
$ws_CRP.Range("H49").Value = 5000
$ws_CRP.Range("I49").Value = 0
$ws_CRP.Range("J49").Value = 5000
$ws_CRP.Range("K49").Value = 0
$ws_CRP.Range("L49").Value = 5000
$ws_CRP.Range("M49").ClearContents()
$ws_CRP.Range("N49").Value = -5364

# Row 107 (Item ID 27689)
$ws_CRP.Range("H107").Value = 884.7778
$ws_CRP.Range("I107").Value = 620.375
$ws_CRP.Range("J107").Value = 3000
$ws_CRP.Range("K107").Value = 620.375
$ws_CRP.Range("L107").Value = 3000
$ws_CRP.Range("M107").Value = 1299.625
$ws_CRP.Range("N107").Value = -6840

# Row 113 (Item ID 27691)
$ws_CRP.Range("H113").Value = 2556
$ws_CRP.Range("I113").Value = 491
$ws_CRP.Range("J113").Value = 3932.6667
$ws_CRP.Range("K113").Value = 491
$ws_CRP.Range("L113").Value = 3932.6667
$ws_CRP.Range("M113").Value = 1679
$ws_CRP.Range("N113").Value = -8272.6667

# --- CUL ---
# Row 4 (Item ID 4650)
$ws_CUL.Range("H4").Value = 143001100
$ws_CUL.Range("I4").Value = 335058.34
$ws_CUL.Range("J4").Value = 250000660
$ws_CUL.Range("K4").Value = 1005175.02
$ws_CUL.Range("L4").Value = 750001980
$ws_CUL.Range("M4").Value = -1005063.02
$ws_CUL.Range("N4").Value = -750002204

# Row 7 (Item ID 4728)
$ws_CUL.Range("H7").Value = 0
$ws_CUL.Range("I7").Value = 0
$ws_CUL.Range("K7").Value = 0
$ws_CUL.Range("M7").ClearContents()

# Row 80 (Item ID 12890)
$ws_CUL.Range("H80").Value = 0
$ws_CUL.Range("I80").Value = 0
$ws_CUL.Range("K80").Value = 0
$ws_CUL.Range("M80").ClearContents()

# Row 83 (Item ID 12890)
$ws_CUL.Range("H83").Value = 0
$ws_CUL.Range("I83").Value = 0
$ws_CUL.Range("K83").Value = 0
$ws_CUL.Range("M83").ClearContents()

# Row 86 (Item ID 12892)
$ws_CUL.Range("H86").Value = 93.2
$ws_CUL.Range("I86").Value = 89.5
$ws_CUL.Range("J86").Value = 95.666664
$ws_CUL.Range("K86").Value = 268.5
$ws_CUL.Range("L86").Value = 286.999992
$ws_CUL.Range("M86").Value = 917.5
$ws_CUL.Range("N86").Value = -2658.999992

# Row 89 (Item ID 12892)
$ws_CUL.Range("H89").Value = 93.2
$ws_CUL.Range("I89").Value = 89.5
$ws_CUL.Range("J89").Value = 95.666664
$ws_CUL.Range("K89").Value = 805.5
$ws_CUL.Range("L89").Value = 860.9999759999999
$ws_CUL.Range("M89").Value = 5122.5
$ws_CUL.Range("N89").Value = -12716.999976

# Row 92 (Item ID 19841)
$ws_CUL.Range("H92").Value = 1022.5
$ws_CUL.Range("I92").Value = 1094.5
$ws_CUL.Range("J92").Value = 950.5
$ws_CUL.Range("K92").Value = 3283.5
$ws_CUL.Range("L92").Value = 2851.5
$ws_CUL.Range("M92").Value = -2035.5
$ws_CUL.Range("N92").Value = -5347.5

# --- GSM ---
# Row 2 (Item ID 5062)
$ws_GSM.Range("H2").Value = 62.333332
$ws_GSM.Range("I2").Value = 35.42857
$ws_GSM.Range("K2").Value = 35.42857
$ws_GSM.Range("M2").Value = 77.57142999999999

# Row 41 (Item ID 2449)
$ws_GSM.Range("H41").Value = 6111.8
$ws_GSM.Range("I41").Value = 3500
$ws_GSM.Range("J41").Value = 10029.5
$ws_GSM.Range("K41").Value = 3500
$ws_GSM.Range("L41").Value = 10029.5
$ws_GSM.Range("M41").Value = -3145
$ws_GSM.Range("N41").Value = -10739.5

# Row 102 (Item ID 36169)
$ws_GSM.Range("H102").Value = 3595.7273
$ws_GSM.Range("I102").Value = 2186.5715
$ws_GSM.Range("J102").Value = 6061.75
$ws_GSM.Range("K102").Value = 2186.5715
$ws_GSM.Range("L102").Value = 6061.75
$ws_GSM.Range("M102").Value = -564.5715
$ws_GSM.Range("N102").Value = -9305.75

# Row 110 (Item ID 25802)
$ws_GSM.Range("H110").Value = 0
$ws_GSM.Range("J110").Value = 0
$ws_GSM.Range("L110").Value = 0
$ws_GSM.Range("N110").ClearContents()

# Row 114 (Item ID 25957)
$ws_GSM.Range("H114").Value = 79000
$ws_GSM.Range("J114").Value = 79000
$ws_GSM.Range("L114").Value = 79000
$ws_GSM.Range("N114").Value = -87678

# --- LTW ---
# Row 2 (Item ID 2631)
$ws_LTW.Range("H2").Value = 3999
$ws_LTW.Range("I2").Value = 0
$ws_LTW.Range("K2").Value = 0
$ws_LTW.Range("M2").ClearContents()

# Row 22 (Item ID 5277)
$ws_LTW.Range("H22").Value = 4814.143
$ws_LTW.Range("I22").Value = 1174.75
$ws_LTW.Range("K22").Value = 1174.75
$ws_LTW.Range("M22").Value = -879.75

# Row 27 (Item ID 5277)
$ws_LTW.Range("H27").Value = 4814.143
$ws_LTW.Range("I27").Value = 1174.75
$ws_LTW.Range("K27").Value = 1174.75
$ws_LTW.Range("M27").Value = -1067.75

# Row 93 (Item ID 19993)
$ws_LTW.Range("H93").Value = 1741
$ws_LTW.Range("I93").Value = 1551.25
$ws_LTW.Range("K93").Value = 1551.25
$ws_LTW.Range("M93").Value = -303.25

# Row 98 (Item ID 18379)
$ws_LTW.Range("H98").Value = 0
$ws_LTW.Range("J98").Value = 0
$ws_LTW.Range("L98").Value = 0
$ws_LTW.Range("N98").ClearContents()

# Row 105 (Item ID 18698)
$ws_LTW.Range("H105").Value = 28327
$ws_LTW.Range("J105").Value = 28327
$ws_LTW.Range("L105").Value = 28327
$ws_LTW.Range("N105").Value = -35315

# --- WVR ---
# Row 2 (Item ID 3307)
$ws_WVR.Range("H2").Value = 15714.286
$ws_WVR.Range("I2").Value = 8000
$ws_WVR.Range("J2").Value = 35000
$ws_WVR.Range("K2").Value = 8000
$ws_WVR.Range("L2").Value = 35000
$ws_WVR.Range("M2").Value = -7888
$ws_WVR.Range("N2").Value = -35224

# Row 32 (Item ID 3066)
$ws_WVR.Range("H32").Value = 6000
$ws_WVR.Range("I32").Value = 0
$ws_WVR.Range("K32").Value = 0
$ws_WVR.Range("M32").ClearContents()

# Row 40 (Item ID 3601)
$ws_WVR.Range("H40").Value = 29999.5
$ws_WVR.Range("I40").Value = 50000
$ws_WVR.Range("J40").Value = 9999
$ws_WVR.Range("K40").Value = 50000
$ws_WVR.Range("L40").Value = 9999
$ws_WVR.Range("M40").Value = -49851
$ws_WVR.Range("N40").Value = -10297

# Row 45 (Item ID 21726)
$ws_WVR.Range("H45").Value = 34711
$ws_WVR.Range("I45").Value = 34711
$ws_WVR.Range("K45").Value = 34711
$ws_WVR.Range("M45").Value = -34220

# Row 105 (Item ID 18710)
$ws_WVR.Range("H105").Value = 20645.25
$ws_WVR.Range("J105").Value = 20645.25
$ws_WVR.Range("L105").Value = 20645.25
$ws_WVR.Range("N105").Value = -27633.25

# Row 107 (Item ID 27746)
$ws_WVR.Range("H107").Value = 1281.7894
$ws_WVR.Range("I107").Value = 1275.6
$ws_WVR.Range("K107").Value = 3826.8
$ws_WVR.Range("M107").Value = -1906.8
